$d = $word.ActiveDocument

# The Title, Author and Date paragraphs were originally split across many
# single-word/space runs (one run per word and one run per separating
# space). Re-asserting each paragraph's full text via Find/Replace collapses
# it back down into a single run per paragraph (matching the target
# document) even though the visible text itself does not change.

[void]$d.Content.Find.Execute(
    "Test 014: MathJax Walker works when switched on", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Test 014: MathJax Walker works when switched on", 2)

[void]$d.Content.Find.Execute(
    "Emma Cliffe, Skills Centre: MASH, University of Bath", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Emma Cliffe, Skills Centre: MASH, University of Bath", 2)

[void]$d.Content.Find.Execute(
    "October 2023", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "October 2023", 2)
